# Apply updated betting-odds values for rows 2, 3, 4, 10 and 13
# on the active worksheet (Sheet1), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 1.62  # G2
$ws.Cells.Item(2, 8).Value = 3.9  # H2
$ws.Cells.Item(2, 9).Value = 5.75  # I2
$ws.Cells.Item(2, 10).Value = 2.2  # J2
$ws.Cells.Item(2, 21).Value = 1.95  # U2
$ws.Cells.Item(2, 22).Value = 1.8  # V2
$ws.Cells.Item(2, 23).Value = 6.5  # W2
$ws.Cells.Item(2, 40).Value = 3.5  # AN2
$ws.Cells.Item(2, 41).Value = 8  # AO2

# Row 3
$ws.Cells.Item(3, 13).Value = 1.13  # M3
$ws.Cells.Item(3, 14).Value = 6  # N3
$ws.Cells.Item(3, 42).Value = 26  # AP3
$ws.Cells.Item(3, 47).Value = 10  # AU3
$ws.Cells.Item(3, 48).Value = 81  # AV3
$ws.Cells.Item(3, 53).Value = 201  # BA3

# Row 4
$ws.Cells.Item(4, 7).Value = 1.7  # G4
$ws.Cells.Item(4, 8).Value = 3.3  # H4
$ws.Cells.Item(4, 9).Value = 5.75  # I4
$ws.Cells.Item(4, 15).Value = 1.4  # O4
$ws.Cells.Item(4, 16).Value = 2.75  # P4
$ws.Cells.Item(4, 17).Value = 2.3  # Q4
$ws.Cells.Item(4, 18).Value = 1.6  # R4
$ws.Cells.Item(4, 21).Value = 2.2  # U4
$ws.Cells.Item(4, 22).Value = 1.62  # V4
$ws.Cells.Item(4, 24).Value = 7  # X4
$ws.Cells.Item(4, 28).Value = 34  # AB4
$ws.Cells.Item(4, 35).Value = 26  # AI4
$ws.Cells.Item(4, 36).Value = 19  # AJ4
$ws.Cells.Item(4, 40).Value = 3.5  # AN4
$ws.Cells.Item(4, 42).Value = 23  # AP4
$ws.Cells.Item(4, 43).Value = 34  # AQ4
$ws.Cells.Item(4, 44).Value = 51  # AR4
$ws.Cells.Item(4, 45).Value = 201  # AS4
$ws.Cells.Item(4, 47).Value = 9.5  # AU4
$ws.Cells.Item(4, 52).Value = 126  # AZ4
$ws.Cells.Item(4, 53).Value = 151  # BA4

# Row 10
$ws.Cells.Item(10, 7).Value = 2.2  # G10
$ws.Cells.Item(10, 9).Value = 3.75  # I10
$ws.Cells.Item(10, 10).Value = 3  # J10
$ws.Cells.Item(10, 12).Value = 4  # L10
$ws.Cells.Item(10, 26).Value = 21  # Z10
$ws.Cells.Item(10, 27).Value = 21  # AA10
$ws.Cells.Item(10, 35).Value = 17  # AI10
$ws.Cells.Item(10, 36).Value = 13  # AJ10
$ws.Cells.Item(10, 41).Value = 13  # AO10
$ws.Cells.Item(10, 54).Value = 251  # BB10
$ws.Cells.Item(10, 56).Value = 126  # BD10

# Row 13
$ws.Cells.Item(13, 7).Value = 1.22  # G13
$ws.Cells.Item(13, 9).Value = 15  # I13
$ws.Cells.Item(13, 11).Value = 2.75  # K13
$ws.Cells.Item(13, 12).Value = 10  # L13
$ws.Cells.Item(13, 17).Value = 1.6  # Q13
$ws.Cells.Item(13, 18).Value = 2.3  # R13
$ws.Cells.Item(13, 26).Value = 7  # Z13
$ws.Cells.Item(13, 30).Value = 11  # AD13
$ws.Cells.Item(13, 31).Value = 26  # AE13
$ws.Cells.Item(13, 34).Value = 29  # AH13
$ws.Cells.Item(13, 36).Value = 34  # AJ13
$ws.Cells.Item(13, 39).Value = 81  # AM13
$ws.Cells.Item(13, 43).Value = 13  # AQ13
$ws.Cells.Item(13, 45).Value = 126  # AS13
$ws.Cells.Item(13, 47).Value = 11  # AU13
$ws.Cells.Item(13, 51).Value = 51  # AY13


Write-Output "Updated odds cells in rows 2, 3, 4, 10, 13."
